# #26: Add setter for shape's paragraph
#
# 1) Bump the "datetimeFigureOut" date field text (12/25/2020 -> 12/27/2020)
#    on every slide layout + the slide master (the field lives on the
#    placeholder shapes there, not on the slides themselves).
# 2) Rework "TextBox 3" on slide 2: split paragraph 1 & 3 into a bold
#    "Test" run followed by the remainder, reword paragraph 2, and let the
#    (no-wrap, autofit) shape grow to fit the new text width.

$p = $ppt.ActivePresentation

function Fix-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Type -eq 14) {
            if ($sh.PlaceholderFormat.Type -eq 16) {
                $sh.TextFrame.TextRange.Text = "12/27/2020"
            }
        }
    }
}

$master = $p.SlideMaster

Fix-DatePlaceholder $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $lay = $master.CustomLayouts.Item($li)
    Fix-DatePlaceholder $lay.Shapes
}

$slide2 = $p.Slides.Item(2)
$shape = $slide2.Shapes.Item(1)
$tr = $shape.TextFrame.TextRange

$para1 = $tr.Paragraphs(1, 1)
$para1.Text = "-p1"
$para1.InsertBefore("Test")
$tr.Characters($para1.Start, 4).Font.Bold = $true

$para2 = $tr.Paragraphs(2, 1)
# Setting the text directly would keep the literal (shared) "Test-p2"
# suffix as its own run; detour through an unrelated placeholder string
# first so the final assignment lands as a single clean run.
$para2.Text = "zzzzzzzzzzz"
$para2 = $tr.Paragraphs(2, 1)
$para2.Text = "id4-Test-p2"

$para3 = $tr.Paragraphs(3, 1)
$para3.Text = "-p3"
$para3.InsertBefore("Test")
$tr.Characters($para3.Start, 4).Font.Bold = $true

$shape.Width = 1516121 / 12700
